$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("inputs")

# Activate the "inputs" sheet (it becomes the selected/active tab)
$ws.Activate()

# Fix the date-start label's spelling (accented "Date début" -> "Date debut")
$ws.Range("A4").Value = "Date debut"

# Format the start/end date cells as dates
$ws.Range("B4").NumberFormat = "mm-dd-yy"
$ws.Range("B4").Copy()
$ws.Range("B5").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Leave the final selection on A5
$ws.Range("A5").Select()
